{"js": "// Replace each arithmetic-problem cell's text in the 20x5 table with its\n// updated expression, addressed by (row, col) so the text collision\n// between the original cell (0,0) \"26+30=\" -> \"63+12=\" and the original\n// cell (15,2) \"63+12=\" -> \"12-4=\" can't cause a mismatch.\nconst cellMap = [\n  [0, 0, \"26+30=\", \"63+12=\"],\n  [0, 1, \"91-49=\", \"86-24=\"],\n  [0, 2, \"19-7=\", \"47-27=\"],\n  [0, 3, \"15+42=\", \"82-49=\"],\n  [0, 4, \"7+61=\", \"3+48=\"],\n  [1, 0, \"35+43=\", \"73-47=\"],\n  [1, 1, \"11+12=\", \"2+7=\"],\n  [1, 2, \"76+10=\", \"54+10=\"],\n  [1, 3, \"25+60=\", \"85-27=\"],\n  [1, 4, \"64+8=\", \"80-76=\"],\n  [2, 0, \"38+46=\", \"17-13=\"],\n  [2, 1, \"90-89=\", \"49+35=\"],\n  [2, 2, \"77+12=\", \"68-63=\"],\n  [2, 3, \"20-0=\", \"62-9=\"],\n  [2, 4, \"12+60=\", \"61-55=\"],\n  [3, 0, \"59-14=\", \"29+0=\"],\n  [3, 1, \"71-17=\", \"68-35=\"],\n  [3, 2, \"30+23=\", \"15-6=\"],\n  [3, 3, \"64+21=\", \"75-57=\"],\n  [3, 4, \"18+52=\", \"66-7=\"],\n  [4, 0, \"28-24=\", \"13+22=\"],\n  [4, 1, \"25-18=\", \"28+71=\"],\n  [4, 2, \"78-33=\", \"11+54=\"],\n  [4, 3, \"53+12=\", \"95+0=\"],\n  [4, 4, \"49+32=\", \"57-12=\"],\n  [5, 0, \"34+50=\", \"24-15=\"],\n  [5, 1, \"82-3=\", \"62+26=\"],\n  [5, 2, \"87-52=\", \"35+22=\"],\n  [5, 3, \"30+0=\", \"79-76=\"],\n  [5, 4, \"2+32=\", \"45+8=\"],\n  [6, 0, \"8-4=\", \"94-6=\"],\n  [6, 1, \"6+79=\", \"89-74=\"],\n  [6, 2, \"46+6=\", \"13+33=\"],\n  [6, 3, \"85+2=\", \"9+8=\"],\n  [6, 4, \"41+57=\", \"25+73=\"],\n  [7, 0, \"42+20=\", \"92-57=\"],\n  [7, 1, \"69+7=\", \"37+25=\"],\n  [7, 2, \"44-20=\", \"10+64=\"],\n  [7, 3, \"11+11=\", \"70-64=\"],\n  [7, 4, \"32+30=\", \"84-23=\"],\n  [8, 0, \"31-28=\", \"67-27=\"],\n  [8, 1, \"62-17=\", \"31+34=\"],\n  [8, 2, \"40+9=\", \"97-4=\"],\n  [8, 3, \"90+7=\", \"2-0=\"],\n  [8, 4, \"97-42=\", \"98-70=\"],\n  [9, 0, \"0+41=\", \"52+18=\"],\n  [9, 1, \"38+8=\", \"60+20=\"],\n  [9, 2, \"2+43=\", \"81-18=\"],\n  [9, 3, \"79+10=\", \"1+35=\"],\n  [9, 4, \"53-27=\", \"99-81=\"],\n  [10, 0, \"64-60=\", \"67-57=\"],\n  [10, 1, \"93-5=\", \"41+36=\"],\n  [10, 2, \"24-1=\", \"51+39=\"],\n  [10, 3, \"24-8=\", \"67-23=\"],\n  [10, 4, \"4+75=\", \"9+78=\"],\n  [11, 0, \"50+10=\", \"76+20=\"],\n  [11, 1, \"83-76=\", \"24+5=\"],\n  [11, 2, \"34+11=\", \"43+21=\"],\n  [11, 3, \"21+34=\", \"47+12=\"],\n  [11, 4, \"80-70=\", \"10+38=\"],\n  [12, 0, \"97-50=\", \"33+2=\"],\n  [12, 1, \"89-28=\", \"91-3=\"],\n  [12, 2, \"93-56=\", \"26+41=\"],\n  [12, 3, \"21-17=\", \"24+44=\"],\n  [12, 4, \"94-87=\", \"57+9=\"],\n  [13, 0, \"50-26=\", \"20+20=\"],\n  [13, 1, \"17-5=\", \"32+54=\"],\n  [13, 2, \"41+3=\", \"80-79=\"],\n  [13, 3, \"71-36=\", \"40-1=\"],\n  [13, 4, \"80-27=\", \"13+10=\"],\n  [14, 0, \"51+26=\", \"44+10=\"],\n  [14, 1, \"65+32=\", \"28+31=\"],\n  [14, 2, \"70-66=\", \"83+1=\"],\n  [14, 3, \"96-53=\", \"16+36=\"],\n  [14, 4, \"64+9=\", \"17+9=\"],\n  [15, 0, \"25+3=\", \"23+0=\"],\n  [15, 1, \"13-5=\", \"80-4=\"],\n  [15, 2, \"63+12=\", \"12-4=\"],\n  [15, 3, \"88+10=\", \"63-11=\"],\n  [15, 4, \"92-86=\", \"19+21=\"],\n  [16, 0, \"28+51=\", \"93-34=\"],\n  [16, 1, \"56+15=\", \"48+30=\"],\n  [16, 2, \"9+81=\", \"12+86=\"],\n  [16, 3, \"0+70=\", \"59-34=\"],\n  [16, 4, \"87-4=\", \"4+57=\"],\n  [17, 0, \"82-61=\", \"66-52=\"],\n  [17, 1, \"45-6=\", \"25+14=\"],\n  [17, 2, \"75+12=\", \"7+63=\"],\n  [17, 3, \"54+43=\", \"74-56=\"],\n  [17, 4, \"34+40=\", \"99-21=\"],\n  [18, 0, \"34+32=\", \"27-19=\"],\n  [18, 1, \"93-18=\", \"45+46=\"],\n  [18, 2, \"9+61=\", \"85-9=\"],\n  [18, 3, \"43+14=\", \"60-1=\"],\n  [18, 4, \"58+36=\", \"34+6=\"],\n  [19, 0, \"7+87=\", \"91-42=\"],\n  [19, 1, \"78-46=\", \"76-43=\"],\n  [19, 2, \"10-3=\", \"17+69=\"],\n  [19, 3, \"13+53=\", \"54+0=\"],\n  [19, 4, \"68-45=\", \"84-4=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document body, found none.\");\n}\nconst table = tables.items[0];\n\n// Resolve every target cell's search hit up front.\nconst hits = [];\nfor (const [row, col, oldText, newText] of cellMap) {\n  const cell = table.getCell(row, col);\n  const results = cell.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  hits.push({ row, col, oldText, newText, results });\n}\nawait context.sync();\n\n// Replace the matched text run in place so paragraph/run formatting\n// (font, size, justification) carries over untouched.\nfor (const hit of hits) {\n  if (hit.results.items.length === 0) {\n    throw new Error(\n      `Cell (${hit.row}, ${hit.col}) did not contain expected text \"${hit.oldText}\".`\n    );\n  }\n  hit.results.items[0].insertText(hit.newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Replace each arithmetic-problem cell text in the 20x5 table with its\n# updated expression, addressed by (row, col) (1-based, matching the Word\n# COM Table.Cell(row, col) indexing). We assign Cell.Range.Text directly\n# (rather than Find/Replace) so each edit is strictly scoped to its own\n# cell -- this matters because the original cell (1,1) \"26+30=\" becomes\n# \"63+12=\", which collides with the original cell (16,3) \"63+12=\" that\n# itself becomes \"12-4=\": a document-wide Find could otherwise re-match\n# the text we just wrote into (1,1) while processing (16,3).\n$cellMap = @(\n  @(1, 1, \"26+30=\", \"63+12=\"),\n  @(1, 2, \"91-49=\", \"86-24=\"),\n  @(1, 3, \"19-7=\", \"47-27=\"),\n  @(1, 4, \"15+42=\", \"82-49=\"),\n  @(1, 5, \"7+61=\", \"3+48=\"),\n  @(2, 1, \"35+43=\", \"73-47=\"),\n  @(2, 2, \"11+12=\", \"2+7=\"),\n  @(2, 3, \"76+10=\", \"54+10=\"),\n  @(2, 4, \"25+60=\", \"85-27=\"),\n  @(2, 5, \"64+8=\", \"80-76=\"),\n  @(3, 1, \"38+46=\", \"17-13=\"),\n  @(3, 2, \"90-89=\", \"49+35=\"),\n  @(3, 3, \"77+12=\", \"68-63=\"),\n  @(3, 4, \"20-0=\", \"62-9=\"),\n  @(3, 5, \"12+60=\", \"61-55=\"),\n  @(4, 1, \"59-14=\", \"29+0=\"),\n  @(4, 2, \"71-17=\", \"68-35=\"),\n  @(4, 3, \"30+23=\", \"15-6=\"),\n  @(4, 4, \"64+21=\", \"75-57=\"),\n  @(4, 5, \"18+52=\", \"66-7=\"),\n  @(5, 1, \"28-24=\", \"13+22=\"),\n  @(5, 2, \"25-18=\", \"28+71=\"),\n  @(5, 3, \"78-33=\", \"11+54=\"),\n  @(5, 4, \"53+12=\", \"95+0=\"),\n  @(5, 5, \"49+32=\", \"57-12=\"),\n  @(6, 1, \"34+50=\", \"24-15=\"),\n  @(6, 2, \"82-3=\", \"62+26=\"),\n  @(6, 3, \"87-52=\", \"35+22=\"),\n  @(6, 4, \"30+0=\", \"79-76=\"),\n  @(6, 5, \"2+32=\", \"45+8=\"),\n  @(7, 1, \"8-4=\", \"94-6=\"),\n  @(7, 2, \"6+79=\", \"89-74=\"),\n  @(7, 3, \"46+6=\", \"13+33=\"),\n  @(7, 4, \"85+2=\", \"9+8=\"),\n  @(7, 5, \"41+57=\", \"25+73=\"),\n  @(8, 1, \"42+20=\", \"92-57=\"),\n  @(8, 2, \"69+7=\", \"37+25=\"),\n  @(8, 3, \"44-20=\", \"10+64=\"),\n  @(8, 4, \"11+11=\", \"70-64=\"),\n  @(8, 5, \"32+30=\", \"84-23=\"),\n  @(9, 1, \"31-28=\", \"67-27=\"),\n  @(9, 2, \"62-17=\", \"31+34=\"),\n  @(9, 3, \"40+9=\", \"97-4=\"),\n  @(9, 4, \"90+7=\", \"2-0=\"),\n  @(9, 5, \"97-42=\", \"98-70=\"),\n  @(10, 1, \"0+41=\", \"52+18=\"),\n  @(10, 2, \"38+8=\", \"60+20=\"),\n  @(10, 3, \"2+43=\", \"81-18=\"),\n  @(10, 4, \"79+10=\", \"1+35=\"),\n  @(10, 5, \"53-27=\", \"99-81=\"),\n  @(11, 1, \"64-60=\", \"67-57=\"),\n  @(11, 2, \"93-5=\", \"41+36=\"),\n  @(11, 3, \"24-1=\", \"51+39=\"),\n  @(11, 4, \"24-8=\", \"67-23=\"),\n  @(11, 5, \"4+75=\", \"9+78=\"),\n  @(12, 1, \"50+10=\", \"76+20=\"),\n  @(12, 2, \"83-76=\", \"24+5=\"),\n  @(12, 3, \"34+11=\", \"43+21=\"),\n  @(12, 4, \"21+34=\", \"47+12=\"),\n  @(12, 5, \"80-70=\", \"10+38=\"),\n  @(13, 1, \"97-50=\", \"33+2=\"),\n  @(13, 2, \"89-28=\", \"91-3=\"),\n  @(13, 3, \"93-56=\", \"26+41=\"),\n  @(13, 4, \"21-17=\", \"24+44=\"),\n  @(13, 5, \"94-87=\", \"57+9=\"),\n  @(14, 1, \"50-26=\", \"20+20=\"),\n  @(14, 2, \"17-5=\", \"32+54=\"),\n  @(14, 3, \"41+3=\", \"80-79=\"),\n  @(14, 4, \"71-36=\", \"40-1=\"),\n  @(14, 5, \"80-27=\", \"13+10=\"),\n  @(15, 1, \"51+26=\", \"44+10=\"),\n  @(15, 2, \"65+32=\", \"28+31=\"),\n  @(15, 3, \"70-66=\", \"83+1=\"),\n  @(15, 4, \"96-53=\", \"16+36=\"),\n  @(15, 5, \"64+9=\", \"17+9=\"),\n  @(16, 1, \"25+3=\", \"23+0=\"),\n  @(16, 2, \"13-5=\", \"80-4=\"),\n  @(16, 3, \"63+12=\", \"12-4=\"),\n  @(16, 4, \"88+10=\", \"63-11=\"),\n  @(16, 5, \"92-86=\", \"19+21=\"),\n  @(17, 1, \"28+51=\", \"93-34=\"),\n  @(17, 2, \"56+15=\", \"48+30=\"),\n  @(17, 3, \"9+81=\", \"12+86=\"),\n  @(17, 4, \"0+70=\", \"59-34=\"),\n  @(17, 5, \"87-4=\", \"4+57=\"),\n  @(18, 1, \"82-61=\", \"66-52=\"),\n  @(18, 2, \"45-6=\", \"25+14=\"),\n  @(18, 3, \"75+12=\", \"7+63=\"),\n  @(18, 4, \"54+43=\", \"74-56=\"),\n  @(18, 5, \"34+40=\", \"99-21=\"),\n  @(19, 1, \"34+32=\", \"27-19=\"),\n  @(19, 2, \"93-18=\", \"45+46=\"),\n  @(19, 3, \"9+61=\", \"85-9=\"),\n  @(19, 4, \"43+14=\", \"60-1=\"),\n  @(19, 5, \"58+36=\", \"34+6=\"),\n  @(20, 1, \"7+87=\", \"91-42=\"),\n  @(20, 2, \"78-46=\", \"76-43=\"),\n  @(20, 3, \"10-3=\", \"17+69=\"),\n  @(20, 4, \"13+53=\", \"54+0=\"),\n  @(20, 5, \"68-45=\", \"84-4=\"),\n)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\nforeach ($entry in $cellMap) {\n    $row = $entry[0]\n    $col = $entry[1]\n    $oldText = $entry[2]\n    $newText = $entry[3]\n\n    $cell = $tbl.Cell($row, $col)\n    $rng = $cell.Range\n\n    $current = $rng.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne $oldText) {\n        Write-Output \"WARNING: cell ($row, $col) was '$current', expected '$oldText'\"\n    }\n\n    $rng.Text = $newText\n}\n\nWrite-Output \"done\"\n"}
